{"js": "// Load all paragraphs in the document body so we can locate the ones we\n// need to edit by their text content.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// 1) Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\nfor (const p of items) {\n  if (p.text === \"September 19, 2025\") {\n    p.insertText(\"September 21, 2025\", \"Replace\");\n    break;\n  }\n}\n\n// 2) Split the single-line mailing address into two separate paragraphs:\n//    \"969 Story Road, San Jose CA 95122\"\n//      -> \"969 Story Road\"\n//      -> \"San Jose, CA 95122\"\nfor (const p of items) {\n  if (p.text === \"969 Story Road, San Jose CA 95122\") {\n    p.insertText(\"969 Story Road\", \"Replace\");\n    p.insertParagraph(\"San Jose, CA 95122\", \"After\");\n    break;\n  }\n}\n\n// 3) Remove the blank \"No Spacing\" paragraph that immediately follows the\n//    \"...Board of Directors\" line.\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Board of Directors\") !== -1) {\n    items[i + 1].delete();\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13)\n    if ($t -eq \"September 19, 2025\") {\n        $rng = $p.Range\n        $rng.MoveEnd(1, -1) | Out-Null\n        $rng.Text = \"September 21, 2025\"\n        break\n    }\n}\n\n# 2) Split the single-line mailing address into two separate paragraphs:\n#    \"969 Story Road, San Jose CA 95122\"\n#      -> \"969 Story Road\"\n#      -> \"San Jose, CA 95122\"\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13)\n    if ($t -eq \"969 Story Road, San Jose CA 95122\") {\n        $rng = $p.Range\n        $rng.MoveEnd(1, -1) | Out-Null\n        $rng.Text = \"969 Story Road\"\n        $p.Range.InsertParagraphAfter() | Out-Null\n        $newPara = $d.Paragraphs.Item($i + 1)\n        $newRng = $newPara.Range\n        $newRng.MoveEnd(1, -1) | Out-Null\n        $newRng.Text = \"San Jose, CA 95122\"\n        break\n    }\n}\n\n# 3) Remove the blank \"No Spacing\" paragraph that immediately follows the\n#    \"...Board of Directors\" line.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13)\n    if ($t -like \"*Board of Directors*\") {\n        $next = $d.Paragraphs.Item($i + 1)\n        $next.Range.Delete()\n        break\n    }\n}\n"}
